$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 186 ---
$ws.Range("D186").Value = 44448
$ws.Range("M186").Value = 456
$ws.Range("N186").Value = 9000
$ws.Range("O186").Value = 9000
$ws.Range("P186").Value = 9000
$ws.Range("R186").Value = "Brasil"
$ws.Range("S186").Value = 2250

# --- Update existing row 187 ---
$ws.Range("D187").Value = 44448
$ws.Range("N187").Value = 9000
$ws.Range("O187").Value = 9000
$ws.Range("P187").Value = 9000
$ws.Range("R187").Value = "Brasil"
$ws.Range("S187").Value = 2250

# --- Update existing row 188 ---
$ws.Range("D188").Value = 44167
$ws.Range("M188").Value = 228
$ws.Range("N188").Value = 7000
$ws.Range("O188").Value = 7000
$ws.Range("P188").Value = 7000
$ws.Range("S188").Value = 1750

# --- Update existing row 189 ---
$ws.Range("D189").Value = 44167
$ws.Range("L189").Value = "Segunda"
$ws.Range("N189").Value = 7000
$ws.Range("O189").Value = 7000
$ws.Range("P189").Value = 7000
$ws.Range("R189").Value = "Perú"
$ws.Range("S189").Value = 1750

# --- Add new row 190 ---
$ws.Range("A190").Value = 3
$ws.Range("B190").Value = "Femacal de La Calera"
$ws.Range("C190").Value = "Coquimbo"
$ws.Range("D190").NumberFormat = $ws.Range("D185").NumberFormat
$ws.Range("D190").Value = 44238
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100108
$ws.Range("H190").Value = "Tropicales y subtropicales"
$ws.Range("I190").Value = 100108002
$ws.Range("J190").Value = "Mango"
$ws.Range("K190").Value = "Sin especificar"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 180
$ws.Range("N190").Value = 6000
$ws.Range("O190").Value = 6000
$ws.Range("P190").Value = 6000
$ws.Range("Q190").Value = "$/bandeja 4 kilos"
$ws.Range("R190").Value = "Perú"
$ws.Range("S190").Value = 1500
$ws.Range("T190").Value = 4

# --- Add new row 191 ---
$ws.Range("A191").Value = 3
$ws.Range("B191").Value = "Femacal de La Calera"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").NumberFormat = $ws.Range("D185").NumberFormat
$ws.Range("D191").Value = 44399
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100108
$ws.Range("H191").Value = "Tropicales y subtropicales"
$ws.Range("I191").Value = 100108002
$ws.Range("J191").Value = "Mango"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 228
$ws.Range("N191").Value = 8000
$ws.Range("O191").Value = 8000
$ws.Range("P191").Value = 8000
$ws.Range("Q191").Value = "$/bandeja 4 kilos"
$ws.Range("R191").Value = "Brasil"
$ws.Range("S191").Value = 2000
$ws.Range("T191").Value = 4
